$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 46)
$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

# Update selection / view to match the final state
$ws.Range("F47").Select()
$excel.ActiveWindow.ScrollRow = 8
